$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "Done"
$ws.Range("D15").Value = "Done"
$ws.Range("D20").Value = "Done"

$ws.Range("D21").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
